$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New period header BB1 (next quarter after BA1 = 45891) - copy formatting from BA1
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# Carry-forward values identical to column BA for rows 3-18
$ws.Range("BB3").Value = -3.560752169208581
$ws.Range("BB4").Value = 1.224484594823672
$ws.Range("BB5").Value = 0.6212498672564903
$ws.Range("BB6").Value = 0.951852872712089
$ws.Range("BB7").Value = -0.3608752035976437
$ws.Range("BB8").Value = 0.09627146709163537
$ws.Range("BB9").Value = 0.1477266864992943
$ws.Range("BB10").Value = -0.4279125887877044
$ws.Range("BB11").Value = -0.002674352087272958
$ws.Range("BB12").Value = 0.3477863758372779
$ws.Range("BB13").Value = -0.8261807291073398
$ws.Range("BB14").Value = -1.099040380746541
$ws.Range("BB15").Value = 1.197694531567151
$ws.Range("BB16").Value = -0.7498286166554458
$ws.Range("BB17").Value = 0.3439499888177044
$ws.Range("BB18").Value = 0.2473045135454655

# Newly re-forecasted values for the most recent rows
$ws.Range("BB19").Value = -2.06674933094535
$ws.Range("BB20").Value = -0.9969640812590996
$ws.Range("BB21").Value = -0.9962562460937296
